# Generate Report for Handback
# This script updates the handback-status workbook with freshly generated
# identifiers / timestamps for the two tracked files:
#   0cd5a046-ba03-4c23-be19-b520af6070d6.md -> 5c10fc32-de83-4758-a2f8-b35ef233937d.md
#   a820592f-cca7-4bad-89b5-a1e24ba01e1c.md -> ffff34e0cd2c-9d21-4091-9966-f0411354bd43.md
# and refreshes the xliff correspond file names / timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "0cd5a046-ba03-4c23-be19-b520af6070d6"
$newGuid1 = "5c10fc32-de83-4758-a2f8-b35ef233937d"
$oldGuid2 = "a820592f-cca7-4bad-89b5-a1e24ba01e1c"
$newGuid2 = "ffff34e0cd2c-9d21-4091-9966-f0411354bd43"

$newXlfZhCn = "$newGuid1.546226525ee30590ae9a50a4ced6cfc33f38cde4.zh-cn.xlf"
$newXlfDeDe = "$newGuid1.546226525ee30590ae9a50a4ced6cfc33f38cde4.de-de.xlf"

$newOverviewDate = "2016-08-25 13:05:41"
$newZhCnHandoffDate = "2016-08-25 13:05:36"
$newZhCnHandbackDate = "2016-08-25 13:06:06"
$newDeDeHandbackDate = "2016-08-25 13:06:17"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
$wsOverview.Range("G2").Value = $newOverviewDate

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("G3").Value = $newOverviewDate

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newGuid1.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\$newGuid2.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid1.md"
$wsZhCn.Range("G2").Value = $newXlfZhCn
$wsZhCn.Range("H2").Value = $newZhCnHandoffDate
$wsZhCn.Range("I2").Value = "$newGuid1.md"
$wsZhCn.Range("J2").Value = $newXlfZhCn
$wsZhCn.Range("K2").Value = $newZhCnHandbackDate

$wsZhCn.Range("A3").Value = "$newGuid2.md"
$wsZhCn.Range("G3").Value = $newXlfZhCn
$wsZhCn.Range("H3").Value = $newZhCnHandoffDate
$wsZhCn.Range("I3").Value = "$newGuid2.md"
$wsZhCn.Range("J3").Value = $newXlfZhCn
$wsZhCn.Range("K3").Value = $newZhCnHandbackDate

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$newGuid2.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "$newGuid2.md"
    }
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid1.md"
$wsDeDe.Range("G2").Value = $newXlfDeDe
$wsDeDe.Range("H2").Value = $newOverviewDate
$wsDeDe.Range("I2").Value = "$newGuid1.md"
$wsDeDe.Range("J2").Value = $newXlfDeDe
$wsDeDe.Range("K2").Value = $newDeDeHandbackDate

$wsDeDe.Range("A3").Value = "$newGuid2.md"
$wsDeDe.Range("G3").Value = $newXlfDeDe
$wsDeDe.Range("H3").Value = $newOverviewDate
$wsDeDe.Range("I3").Value = "$newGuid2.md"
$wsDeDe.Range("J3").Value = $newXlfDeDe
$wsDeDe.Range("K3").Value = $newDeDeHandbackDate

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$newGuid2.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "$newGuid2.md"
    }
}
